$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price/volume snapshot refresh.
# Note: some new "Price" values are plain decimals (e.g. "11.80").
# These cells are text (t="inlineStr") in the workbook, so a leading
# apostrophe is used to force Excel to keep them as text instead of
# auto-converting to a number and silently dropping trailing zeros.

$ws.Range("D2").Value = "68.364.73"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "2.711.88"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'609.34"
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("D6").Value = "'166.57"
$ws.Range("E6").Value = "  +4.83%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +2.94%  "
$ws.Range("D9").Value = "2.711.26"
$ws.Range("E9").Value = "  +2.38%  "
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "'28.38"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").Value = "3.205.52"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").Value = "'0.0000188"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "68.293.27"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "2.704.15"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").Value = "'11.80"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("D20").Value = "'369.92"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").Value = "'7.62"
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").Value = "'4.48"
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("D23").Value = "'4.92"
$ws.Range("E23").Value = "  +2.81%  "
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "'73.01"
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "2.846.61"
$ws.Range("E28").Value = "  +1.67%  "
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "'578.79"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("D32").Value = "'8.10"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").Value = "'1.98"
$ws.Range("E34").Value = "  +5.58%  "
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  -3.03%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'19.86"
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'160.57"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "'5.38"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.86"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").Value = "'17.97"
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("D44").Value = "'2.60"
$ws.Range("E44").Value = "  -1.98%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "0.0₆0310"
$ws.Range("E46").Value = "  -3.56%  "
$ws.Range("D47").Value = "'40.80"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("D48").Value = "'0.596"
$ws.Range("E48").Value = "  +3.53%  "
$ws.Range("D49").Value = "'154.78"
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("D50").Value = "'3.90"
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("E51").Value = "  +3.57%  "
